# Auto-generated edit script updating Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2794.889
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2894.25
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2894.25
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -3244.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7077.6665
$ws.Range("I51").Value = 7566.3335
$ws.Range("J51").Value = 6833.3335
$ws.Range("K51").Value = 7566.3335
$ws.Range("L51").Value = 6833.3335
$ws.Range("M51").Value = -7082.3335
$ws.Range("N51").Value = -7801.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 26263.5
$ws.Range("I62").Value = 5206.769
$ws.Range("K62").Value = 5206.769
$ws.Range("M62").Value = -4582.769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 26263.5
$ws.Range("I65").Value = 5206.769
$ws.Range("K65").Value = 26033.845
$ws.Range("M65").Value = -22913.845

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2267.1333
$ws.Range("I106").Value = 2141.2856
$ws.Range("K106").Value = 2141.2856
$ws.Range("M106").Value = -1510.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2922.8235
$ws.Range("I137").Value = 4373.5
$ws.Range("J137").Value = 1633.3334
$ws.Range("K137").Value = 13120.5
$ws.Range("L137").Value = 4900.0002
$ws.Range("M137").Value = -10570.5
$ws.Range("N137").Value = -10000.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 15698
$ws.Range("J119").Value = 15698
$ws.Range("L119").Value = 15698
$ws.Range("N119").Value = -25374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 15404.143
$ws.Range("J124").Value = 15404.143
$ws.Range("L124").Value = 15404.143
$ws.Range("N124").Value = -25224.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6071.6577
$ws.Range("I134").Value = 2531.45
$ws.Range("J134").Value = 10005.223
$ws.Range("K134").Value = 7594.349999999999
$ws.Range("L134").Value = 30015.669
$ws.Range("M134").Value = -5059.349999999999
$ws.Range("N134").Value = -35085.669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3200.138
$ws.Range("I31").Value = 2405.7778
$ws.Range("J31").Value = 4500
$ws.Range("K31").Value = 2405.7778
$ws.Range("L31").Value = 4500
$ws.Range("M31").Value = -2110.7778
$ws.Range("N31").Value = -5090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3200.138
$ws.Range("I34").Value = 2405.7778
$ws.Range("J34").Value = 4500
$ws.Range("K34").Value = 2405.7778
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -2203.7778
$ws.Range("N34").Value = -4904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3563.8823
$ws.Range("I99").Value = 2597.3
$ws.Range("K99").Value = 2597.3
$ws.Range("M99").Value = -1099.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3563.8823
$ws.Range("I126").Value = 2597.3
$ws.Range("K126").Value = 7791.900000000001
$ws.Range("M126").Value = -5321.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 29593.334
$ws.Range("J130").Value = 29593.334
$ws.Range("L130").Value = 29593.334
$ws.Range("N130").Value = -39633.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3231.1667
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 3777.4
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 11332.2
$ws.Range("M94").Value = -824
$ws.Range("N94").Value = -12684.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4900
$ws.Range("J96").Value = 4900
$ws.Range("L96").Value = 14700
$ws.Range("N96").Value = -18818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1596.2222
$ws.Range("I97").Value = 1372.4
$ws.Range("J97").Value = 1876
$ws.Range("K97").Value = 4117.200000000001
$ws.Range("L97").Value = 5628
$ws.Range("M97").Value = -3621.200000000001
$ws.Range("N97").Value = -6620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1243.7142
$ws.Range("I98").Value = 363.125
$ws.Range("J98").Value = 2417.8333
$ws.Range("K98").Value = 1089.375
$ws.Range("L98").Value = 7253.499899999999
$ws.Range("M98").Value = 408.625
$ws.Range("N98").Value = -10249.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 8444.444
$ws.Range("J101").Value = 8444.444
$ws.Range("L101").Value = 25333.332
$ws.Range("N101").Value = -30201.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 10999.333
$ws.Range("J105").Value = 10999.333
$ws.Range("L105").Value = 32997.999
$ws.Range("N105").Value = -38239.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 38461924
$ws.Range("I107").Value = 125000220
$ws.Range("J107").Value = 454.44446
$ws.Range("K107").Value = 375000660
$ws.Range("L107").Value = 1363.33338
$ws.Range("M107").Value = -374998740
$ws.Range("N107").Value = -5203.33338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3611.111
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 3937.5
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 11812.5
$ws.Range("M110").Value = 1090
$ws.Range("N110").Value = -19992.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 55162.26
$ws.Range("I129").Value = 584
$ws.Range("J129").Value = 74654.5
$ws.Range("K129").Value = 1752
$ws.Range("L129").Value = 223963.5
$ws.Range("M129").Value = 3248
$ws.Range("N129").Value = -233963.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 15285.714
$ws.Range("J63").Value = 15285.714
$ws.Range("L63").Value = 15285.714
$ws.Range("N63").Value = -16657.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 15285.714
$ws.Range("J66").Value = 15285.714
$ws.Range("L66").Value = 45857.142
$ws.Range("N66").Value = -52721.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 25100.5
$ws.Range("J69").Value = 25100.5
$ws.Range("L69").Value = 25100.5
$ws.Range("N69").Value = -26598.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 25100.5
$ws.Range("J72").Value = 25100.5
$ws.Range("L72").Value = 75301.5
$ws.Range("N72").Value = -82789.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 17333.334
$ws.Range("J74").Value = 17333.334
$ws.Range("L74").Value = 17333.334
$ws.Range("N74").Value = -19205.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 17333.334
$ws.Range("J77").Value = 17333.334
$ws.Range("L77").Value = 52000.00199999999
$ws.Range("N77").Value = -61360.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 60000
$ws.Range("J82").Value = 60000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60766

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 60000
$ws.Range("J85").Value = 60000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 40195
$ws.Range("J88").Value = 40195
$ws.Range("L88").Value = 40195
$ws.Range("N88").Value = -41097

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 40195
$ws.Range("J91").Value = 40195
$ws.Range("L91").Value = 40195
$ws.Range("N91").Value = -43315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 33335834
$ws.Range("I136").Value = 55556580
$ws.Range("J136").Value = 4717.5
$ws.Range("K136").Value = 166669740
$ws.Range("L136").Value = 14152.5
$ws.Range("M136").Value = -166667190
$ws.Range("N136").Value = -19252.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 5933.3335
$ws.Range("J47").Value = 5933.3335
$ws.Range("L47").Value = 5933.3335
$ws.Range("N47").Value = -7077.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6580505.5
$ws.Range("I136").Value = 17858250
$ws.Range("J136").Value = 1821.875
$ws.Range("K136").Value = 53574750
$ws.Range("L136").Value = 5465.625
$ws.Range("M136").Value = -53572200
$ws.Range("N136").Value = -10565.625
